$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "column 1"
$ws.Range("A2").Value = "Added to excel"

$ws.Range("A3").Select()
